$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data in column B
$ws.Range("B1").Value = 19
$ws.Range("B2").Value = 23

# A3 becomes a formula (still evaluates to FALSE); B3 sums the new column
$ws.Range("A3").Formula = "=FALSE()"
$ws.Range("B3").Formula = "=SUM(B1:B2)"

# Drop the sheet's explicit (all-default) column-width overrides now that
# the used range no longer extends past column B, so the old blanket
# <cols> entry collapses down to just the columns actually in use.
$tail = $ws.Range($ws.Cells.Item(1, 3), $ws.Cells.Item(1, 16384))
$tail.EntireColumn.Delete()

$ws.Range("B3").Select()
